$wb = $excel.ActiveWorkbook

$wsInput    = $wb.Worksheets.Item("Input")
$wsSummary  = $wb.Worksheets.Item("Summary")
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsTrans    = $wb.Worksheets.Item("Transactions")

# --- Summary sheet: add a new (blank, styleless) column G next to the
#     existing data by copying an already-styleless cell into G2. This both
#     creates the bare <c r="G2"/> stub and extends the sheet's used range
#     (dimension + row spans) from F to G.
$wsTrans.Range("K3").Copy($wsSummary.Range("G2"))

# --- Repayment schedule sheet: insert a new column O (mirroring column N,
#     which is all zero/blank) between the existing N and P columns.
$wsSchedule.Range("N2:N8").Copy($wsSchedule.Range("O2:O8"))

# --- Transactions sheet: correct the running transaction IDs.
$wsTrans.Range("A2").Value2 = 89
$wsTrans.Range("A3").Value2 = 88

# --- Selections (recorded per-sheet cursor position / highlighted range). ---
$wsSummary.Activate()
$wsSummary.Range("B5").Select()

$wsSchedule.Activate()
$wsSchedule.Rows.Item(9).Select()

$wsTrans.Activate()
$wsTrans.Range("A2:L3").Select()

Write-Host "done"
